# Applies the diff between the pre-edit and post-edit versions of
# 2.7b_idc4_npe0_results.xlsx.
#
# Sheet order / names (tab order, 1-based):
#   1 pcroprep   2 pfoodrep   3 pliverep   4 pdietrep
#   5 pradar     6 plandrep  7 plaborrep

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) pcroprep ("pcroprep") - numeric updates on rows 17 and 21
# ---------------------------------------------------------------------
$pcroprep = $wb.Worksheets.Item("pcroprep")

$pcroprep.Range("C17").Value = 2805550240019.7905
$pcroprep.Range("D17").Value = 10276740805933.297
$pcroprep.Range("G17").Value = 225.33376760969122
$pcroprep.Range("H17").Value = 61.516118557445715
$pcroprep.Range("I17").Value = 55.633767609691233

$pcroprep.Range("C21").Value = 82406734495546.547
$pcroprep.Range("D21").Value = 66206293629486.211
$pcroprep.Range("G21").Value = 1451.6775176807409
$pcroprep.Range("H21").Value = 1806.8977617468131
$pcroprep.Range("I21").Value = 645.07751768074093

# ---------------------------------------------------------------------
# 2) pdietrep - numeric updates on rows 2-9 (columns D, E, G)
# ---------------------------------------------------------------------
$pdietrep = $wb.Worksheets.Item("pdietrep")

$pdietrep.Range("D2").Value = 2074.3096540334068
$pdietrep.Range("E2").Value = 96.131410669748377
$pdietrep.Range("G2").Value = -83.475860172276043

$pdietrep.Range("D3").Value = 49.882731629334181
$pdietrep.Range("E3").Value = 71.130783506574417
$pdietrep.Range("G3").Value = -20.245459247016377

$pdietrep.Range("D4").Value = 14.476851675698938
$pdietrep.Range("E4").Value = 22.363748358971495
$pdietrep.Range("G4").Value = -50.256713750471533

$pdietrep.Range("D5").Value = 415.44633040783594
$pdietrep.Range("E5").Value = 128.35576958344996
$pdietrep.Range("G5").Value = 91.778503276983599

$pdietrep.Range("D6").Value = 1622.419296516086
$pdietrep.Range("E6").Value = 75.189090196173922
$pdietrep.Range("G6").Value = -535.36621768959685

$pdietrep.Range("D7").Value = 56.456277979493251
$pdietrep.Range("E7").Value = 80.504398122912491
$pdietrep.Range("G7").Value = -13.671912896857307

$pdietrep.Range("D8").Value = 17.181657184106431
$pdietrep.Range("E8").Value = 26.542114698907401
$pdietrep.Range("G8").Value = -47.551908242064044

$pdietrep.Range("D9").Value = 296.97302395876653
$pdietrep.Range("E9").Value = 91.75240758133998
$pdietrep.Range("G9").Value = -26.694803172085813

# ---------------------------------------------------------------------
# 3) pradar - numeric updates on row 7
# ---------------------------------------------------------------------
$pradar = $wb.Worksheets.Item("pradar")

$pradar.Range("C7").Value = 10276740805933.297
$pradar.Range("E7").Value = 225.33376760969122
$pradar.Range("F7").Value = -16.066232390308755
$pradar.Range("G7").Value = 93.344559904594547

# ---------------------------------------------------------------------
# 4) plandrep - cell removals + numeric updates
# ---------------------------------------------------------------------
$plandrep = $wb.Worksheets.Item("plandrep")

# Row 2
$plandrep.Range("I2").ClearContents()
$plandrep.Range("T2").Value = 195.36637075410198
$plandrep.Range("AA2").Value = 826.26637075410201

# Row 3
$plandrep.Range("T3").Value = 195.36637075410198
$plandrep.Range("AA3").Value = 4221.8143707541012

# Row 5
$plandrep.Range("I5").ClearContents()
$plandrep.Range("T5").Value = 424.23274150820396

# Row 6
$plandrep.Range("I6").ClearContents()
$plandrep.Range("AA6").Value = 1048.8999999999999

# Row 9
$plandrep.Range("I9").ClearContents()

# Rows 12-13: drop the string values but keep the (styled) empty cells
$plandrep.Range("U12:V12").ClearContents()
$plandrep.Range("U13:V13").ClearContents()

# Rows 14-17: M/P values updated (404.0642752505355 -> 424.23274150820396)
foreach ($r in 14..17) {
    $plandrep.Range("M$r").Value = 424.23274150820396
    $plandrep.Range("P$r").Value = 424.23274150820396
}

# Row 18: M/P updated + U/V cleared
$plandrep.Range("M18").Value = 395.51637075410196
$plandrep.Range("P18").Value = 1026.416370754102
$plandrep.Range("U18:V18").ClearContents()

# Rows 19-22: only U/V cleared
foreach ($r in 19..22) {
    $plandrep.Range("U$r" + ":V$r").ClearContents()
}

# Row 23: M/P updated + U/V cleared
$plandrep.Range("M23").Value = 424.23274150820396
$plandrep.Range("P23").Value = 475.43274150820395
$plandrep.Range("U23:V23").ClearContents()

# Rows 24-25: M/P updated
foreach ($r in 24..25) {
    $plandrep.Range("M$r").Value = 424.23274150820396
    $plandrep.Range("P$r").Value = 424.23274150820396
}

# Rows 30-34: only U/V cleared
foreach ($r in 30..34) {
    $plandrep.Range("U$r" + ":V$r").ClearContents()
}

# Row 35: only U/V cleared
$plandrep.Range("U35:V35").ClearContents()

# ---------------------------------------------------------------------
# 5) plaborrep - cell removals + numeric updates
# ---------------------------------------------------------------------
$plaborrep = $wb.Worksheets.Item("plaborrep")

$plaborrep.Range("H2").ClearContents()
$plaborrep.Range("S2").Value = 0.69149936865837103
$plaborrep.Range("AG2").Value = 3.3330521843291856

$plaborrep.Range("H3").ClearContents()
$plaborrep.Range("AG3").Value = 3.2373330074954088

$plaborrep.Range("F7").Value = 0.69149936865837103
$plaborrep.Range("K7").Value = 3.3330521843291856

# ---------------------------------------------------------------------
# 6) Active sheet: plaborrep becomes the selected tab (was pcroprep)
# ---------------------------------------------------------------------
$plaborrep.Select()
